$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 40
$ws.Range("H40").Value = 4233
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4233
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4233
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4583
# row 43
$ws.Range("H43").Value = 4131
$ws.Range("I43").Value = 4033
$ws.Range("J43").Value = 4189.8
$ws.Range("K43").Value = 4033
$ws.Range("L43").Value = 4189.8
$ws.Range("M43").Value = -3964
$ws.Range("N43").Value = -4327.8
# row 58
$ws.Range("H58").Value = 103
$ws.Range("I58").Value = 103
$ws.Range("K58").Value = 309
$ws.Range("M58").Value = -159
# row 101
$ws.Range("H101").Value = 306.375
$ws.Range("I101").Value = 300.2857
$ws.Range("K101").Value = 900.8571000000001
$ws.Range("M101").Value = 721.1428999999999
# row 125
$ws.Range("H125").Value = 5106.4
$ws.Range("I125").Value = 1766
$ws.Range("K125").Value = 15894
$ws.Range("M125").Value = -13434
# row 132
$ws.Range("H132").Value = 9925.166999999999
$ws.Range("I132").Value = 3440.8125
$ws.Range("J132").Value = 30675.1
$ws.Range("K132").Value = 10322.4375
$ws.Range("L132").Value = 92025.29999999999
$ws.Range("M132").Value = -7792.4375
$ws.Range("N132").Value = -97085.29999999999
# row 138
$ws.Range("H138").Value = 224676.42
$ws.Range("I138").Value = 1062.8334
$ws.Range("J138").Value = 305990.47
$ws.Range("K138").Value = 3188.5002
$ws.Range("L138").Value = 917971.4099999999
$ws.Range("M138").Value = 1951.4998
$ws.Range("N138").Value = -928251.4099999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 2129.6924
$ws.Range("I2").Value = 2295.6667
$ws.Range("J2").Value = 1756.25
$ws.Range("K2").Value = 2295.6667
$ws.Range("L2").Value = 1756.25
$ws.Range("M2").Value = -2182.6667
$ws.Range("N2").Value = -1982.25
# row 45
$ws.Range("H45").Value = 3309.6667
$ws.Range("I45").Value = 2916.889
$ws.Range("K45").Value = 2916.889
$ws.Range("M45").Value = -2539.889
# row 63
$ws.Range("H63").Value = 2895.6667
$ws.Range("I63").Value = 2343.5
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 2343.5
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1657.5
$ws.Range("N63").Value = -5372
# row 66
$ws.Range("H66").Value = 2895.6667
$ws.Range("I66").Value = 2343.5
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 11717.5
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -8285.5
$ws.Range("N66").Value = -26864
# row 74
$ws.Range("H74").Value = 2274.037
$ws.Range("I74").Value = 1059.8572
$ws.Range("J74").Value = 6523.6665
$ws.Range("K74").Value = 1059.8572
$ws.Range("L74").Value = 6523.6665
$ws.Range("M74").Value = -185.8571999999999
$ws.Range("N74").Value = -8271.666499999999
# row 77
$ws.Range("H77").Value = 2274.037
$ws.Range("I77").Value = 1059.8572
$ws.Range("J77").Value = 6523.6665
$ws.Range("K77").Value = 5299.286
$ws.Range("L77").Value = 32618.3325
$ws.Range("M77").Value = -931.2860000000001
$ws.Range("N77").Value = -41354.3325
# row 107
$ws.Range("H107").Value = 91122.3
$ws.Range("J107").Value = 91122.3
$ws.Range("L107").Value = 91122.3
$ws.Range("N107").Value = -98802.3
# row 109
$ws.Range("H109").Value = 90000
$ws.Range("J109").Value = 90000
$ws.Range("L109").Value = 90000
$ws.Range("N109").Value = -92774
# row 110
$ws.Range("H110").Value = 1675.3572
$ws.Range("I110").Value = 1686.8
$ws.Range("K110").Value = 1686.8
$ws.Range("M110").Value = 358.2
# row 116
$ws.Range("H116").Value = 2129.6924
$ws.Range("I116").Value = 2295.6667
$ws.Range("J116").Value = 1756.25
$ws.Range("K116").Value = 2295.6667
$ws.Range("L116").Value = 1756.25
$ws.Range("M116").Value = -1.666700000000219
$ws.Range("N116").Value = -6344.25
# row 122
$ws.Range("H122").Value = 3702.3076
$ws.Range("I122").Value = 3204.0334
$ws.Range("K122").Value = 9612.100199999999
$ws.Range("M122").Value = -7162.100199999999
# row 132
$ws.Range("H132").Value = 3905.4375
$ws.Range("I132").Value = 3868.138
$ws.Range("J132").Value = 4266
$ws.Range("K132").Value = 11604.414
$ws.Range("L132").Value = 12798
$ws.Range("M132").Value = -9074.414000000001
$ws.Range("N132").Value = -17858
# row 138
$ws.Range("H138").Value = 89998.664
$ws.Range("J138").Value = 89998.664
$ws.Range("L138").Value = 89998.664
$ws.Range("N138").Value = -100278.664

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 2129.6924
$ws.Range("I3").Value = 2295.6667
$ws.Range("J3").Value = 1756.25
$ws.Range("K3").Value = 2295.6667
$ws.Range("L3").Value = 1756.25
$ws.Range("M3").Value = -2181.6667
$ws.Range("N3").Value = -1984.25
# row 82
$ws.Range("H82").Value = 23220.766
$ws.Range("J82").Value = 32864.184
$ws.Range("L82").Value = 32864.184
$ws.Range("N82").Value = -33630.184
# row 85
$ws.Range("H85").Value = 23220.766
$ws.Range("J85").Value = 32864.184
$ws.Range("L85").Value = 32864.184
$ws.Range("N85").Value = -35516.184
# row 86
$ws.Range("H86").Value = 1952.5526
$ws.Range("I86").Value = 2056.2173
$ws.Range("J86").Value = 1793.6
$ws.Range("K86").Value = 2056.2173
$ws.Range("L86").Value = 1793.6
$ws.Range("M86").Value = -933.2172999999998
$ws.Range("N86").Value = -4039.6
# row 89
$ws.Range("H89").Value = 1952.5526
$ws.Range("I89").Value = 2056.2173
$ws.Range("J89").Value = 1793.6
$ws.Range("K89").Value = 10281.0865
$ws.Range("L89").Value = 8968
$ws.Range("M89").Value = -4665.086499999999
$ws.Range("N89").Value = -20200
# row 134
$ws.Range("H134").Value = 7159.5103
$ws.Range("I134").Value = 3283.353
$ws.Range("K134").Value = 9850.059000000001
$ws.Range("M134").Value = -7315.059000000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2716.3333
$ws.Range("I31").Value = 2002.3158
$ws.Range("J31").Value = 9499.5
$ws.Range("K31").Value = 2002.3158
$ws.Range("L31").Value = 9499.5
$ws.Range("M31").Value = -1707.3158
$ws.Range("N31").Value = -10089.5
# row 34
$ws.Range("H34").Value = 2716.3333
$ws.Range("I34").Value = 2002.3158
$ws.Range("J34").Value = 9499.5
$ws.Range("K34").Value = 2002.3158
$ws.Range("L34").Value = 9499.5
$ws.Range("M34").Value = -1800.3158
$ws.Range("N34").Value = -9903.5
# row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# row 132
$ws.Range("H132").Value = 1144530.4
$ws.Range("I132").Value = 1430024.4
$ws.Range("J132").Value = 2554.5715
$ws.Range("K132").Value = 4290073.199999999
$ws.Range("L132").Value = 7663.7145
$ws.Range("M132").Value = -4287543.199999999
$ws.Range("N132").Value = -12723.7145
# row 134
$ws.Range("H134").Value = 3992.3044
$ws.Range("I134").Value = 1755.6923
$ws.Range("K134").Value = 5267.0769
$ws.Range("M134").Value = -2732.0769

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 4860
$ws.Range("I5").Value = 1065.4615
$ws.Range("J5").Value = 8148.6
$ws.Range("K5").Value = 3196.3845
$ws.Range("L5").Value = 24445.8
$ws.Range("M5").Value = -3084.3845
$ws.Range("N5").Value = -24669.8
# row 9
$ws.Range("H9").Value = 400420
$ws.Range("I9").Value = 400420
$ws.Range("K9").Value = 1201260
$ws.Range("M9").Value = -1201036
# row 63
$ws.Range("H63").Value = 116799.78
$ws.Range("I63").Value = 335399.34
$ws.Range("J63").Value = 7500
$ws.Range("K63").Value = 1006198.02
$ws.Range("L63").Value = 22500
$ws.Range("M63").Value = -1005449.02
$ws.Range("N63").Value = -23998
# row 64
$ws.Range("H64").Value = 2661.3333
$ws.Range("J64").Value = 3442
$ws.Range("L64").Value = 10326
$ws.Range("N64").Value = -10866
# row 66
$ws.Range("H66").Value = 116799.78
$ws.Range("I66").Value = 335399.34
$ws.Range("J66").Value = 7500
$ws.Range("K66").Value = 3018594.06
$ws.Range("L66").Value = 67500
$ws.Range("M66").Value = -3014850.06
$ws.Range("N66").Value = -74988
# row 67
$ws.Range("H67").Value = 2661.3333
$ws.Range("J67").Value = 3442
$ws.Range("L67").Value = 10326
$ws.Range("N67").Value = -12198
# row 113
$ws.Range("H113").Value = 1519.7646
$ws.Range("I113").Value = 847
$ws.Range("J113").Value = 1561.8125
$ws.Range("K113").Value = 2541
$ws.Range("L113").Value = 4685.4375
$ws.Range("M113").Value = -371
$ws.Range("N113").Value = -9025.4375
# row 132
$ws.Range("H132").Value = 2490.439
$ws.Range("I132").Value = 911.9
$ws.Range("J132").Value = 2999.6453
$ws.Range("K132").Value = 8207.1
$ws.Range("L132").Value = 26996.8077
$ws.Range("M132").Value = -5677.1
$ws.Range("N132").Value = -32056.8077
# row 135
$ws.Range("H135").Value = 4860
$ws.Range("I135").Value = 1065.4615
$ws.Range("J135").Value = 8148.6
$ws.Range("K135").Value = 9589.153499999999
$ws.Range("L135").Value = 73337.40000000001
$ws.Range("M135").Value = -7054.153499999999
$ws.Range("N135").Value = -78407.40000000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Range("H122").Value = 3023.08
$ws.Range("I122").Value = 1994
$ws.Range("K122").Value = 5982
$ws.Range("M122").Value = -3532
# row 126
$ws.Range("H126").Value = 3896.4443
$ws.Range("I126").Value = 3828
$ws.Range("K126").Value = 11484
$ws.Range("M126").Value = -9014
# row 132
$ws.Range("H132").Value = 12350446
$ws.Range("I132").Value = 14497398
$ws.Range("K132").Value = 43492194
$ws.Range("M132").Value = -43489664

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 2944.5881
$ws.Range("I132").Value = 2902.2886
$ws.Range("K132").Value = 8706.8658
$ws.Range("M132").Value = -6176.8658
# row 133
$ws.Range("H133").Value = 74708.664
$ws.Range("J133").Value = 74708.664
$ws.Range("L133").Value = 74708.664
$ws.Range("N133").Value = -79768.664
# row 140
$ws.Range("H140").Value = 159261.5
$ws.Range("J140").Value = 159261.5
$ws.Range("L140").Value = 159261.5
$ws.Range("N140").Value = -169621.5
